$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1361.6666
$ws.Range("I17").Value = 750
$ws.Range("J17").Value = 1388.2609
$ws.Range("K17").Value = 2250
$ws.Range("L17").Value = 4164.7827
$ws.Range("M17").Value = -2082
$ws.Range("N17").Value = -4500.7827
$ws.Range("H28").Value = 1110.4166
$ws.Range("I28").Value = 832.6
$ws.Range("K28").Value = 832.6
$ws.Range("M28").Value = -347.6
$ws.Range("H55").Value = 1588.5555
$ws.Range("I55").Value = 966.3333
$ws.Range("K55").Value = 966.3333
$ws.Range("M55").Value = -752.3333
$ws.Range("H62").Value = 4649.2
$ws.Range("I62").Value = 4610.222
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 4610.222
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -3986.222
$ws.Range("N62").Value = -6248
$ws.Range("H65").Value = 4649.2
$ws.Range("I65").Value = 4610.222
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 23051.11
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -19931.11
$ws.Range("N65").Value = -31240
$ws.Range("H74").Value = 10000
$ws.Range("I74").Value = 10000
$ws.Range("K74").Value = 10000
$ws.Range("M74").Value = -9064
$ws.Range("H77").Value = 10000
$ws.Range("I77").Value = 10000
$ws.Range("K77").Value = 50000
$ws.Range("M77").Value = -45320
$ws.Range("H132").Value = 1654.6522
$ws.Range("I132").Value = 1656.4
$ws.Range("K132").Value = 4969.200000000001
$ws.Range("M132").Value = -2439.200000000001
$ws.Range("H137").Value = 2003.0667
$ws.Range("I137").Value = 2022.4546
$ws.Range("K137").Value = 6067.3638
$ws.Range("M137").Value = -3517.3638
$ws.Range("H138").Value = 14614.462
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 14614.462
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 43843.386
$ws.Range("M138").Value = ""
$ws.Range("N138").Value = -54123.386
$ws.Range("H141").Value = 3533.7827
$ws.Range("I141").Value = 2009.2
$ws.Range("J141").Value = 13697.667
$ws.Range("K141").Value = 6027.6
$ws.Range("L141").Value = 41093.001
$ws.Range("M141").Value = -847.6000000000004
$ws.Range("N141").Value = -51453.001
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3493.75
$ws.Range("I2").Value = 2089
$ws.Range("J2").Value = 4898.5
$ws.Range("K2").Value = 2089
$ws.Range("L2").Value = 4898.5
$ws.Range("M2").Value = -1976
$ws.Range("N2").Value = -5124.5
$ws.Range("H97").Value = 3299
$ws.Range("I97").Value = 2732
$ws.Range("K97").Value = 2732
$ws.Range("M97").Value = -2236
$ws.Range("H116").Value = 3493.75
$ws.Range("I116").Value = 2089
$ws.Range("J116").Value = 4898.5
$ws.Range("K116").Value = 2089
$ws.Range("L116").Value = 4898.5
$ws.Range("M116").Value = 205
$ws.Range("N116").Value = -9486.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3493.75
$ws.Range("I3").Value = 2089
$ws.Range("J3").Value = 4898.5
$ws.Range("K3").Value = 2089
$ws.Range("L3").Value = 4898.5
$ws.Range("M3").Value = -1975
$ws.Range("N3").Value = -5126.5
$ws.Range("H86").Value = 7602.3335
$ws.Range("I86").Value = 4000
$ws.Range("J86").Value = 9403.5
$ws.Range("K86").Value = 4000
$ws.Range("L86").Value = 9403.5
$ws.Range("M86").Value = -2877
$ws.Range("N86").Value = -11649.5
$ws.Range("H89").Value = 7602.3335
$ws.Range("I89").Value = 4000
$ws.Range("J89").Value = 9403.5
$ws.Range("K89").Value = 20000
$ws.Range("L89").Value = 47017.5
$ws.Range("M89").Value = -14384
$ws.Range("N89").Value = -58249.5
$ws.Range("H134").Value = 1782.619
$ws.Range("J134").Value = 2674.5
$ws.Range("L134").Value = 8023.5
$ws.Range("N134").Value = -13093.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 484.42856
$ws.Range("I16").Value = 484.42856
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 484.42856
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -197.42856
$ws.Range("N16").Value = ""
$ws.Range("H113").Value = 484.42856
$ws.Range("I113").Value = 484.42856
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 484.42856
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1685.57144
$ws.Range("N113").Value = ""
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 192
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").Value = ""
$ws.Range("H122").Value = 297
$ws.Range("I122").Value = 297
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2673
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -223
$ws.Range("N122").Value = ""
$ws.Range("H126").Value = 3499.75
$ws.Range("I126").Value = 3499.75
$ws.Range("K126").Value = 10499.25
$ws.Range("M126").Value = -5559.25
$ws.Range("H128").Value = 3979888.5
$ws.Range("I128").Value = 3979888.5
$ws.Range("K128").Value = 11939665.5
$ws.Range("M128").Value = -11934685.5
$ws.Range("H129").Value = 2563.9167
$ws.Range("I129").Value = 1593
$ws.Range("J129").Value = 2887.5557
$ws.Range("K129").Value = 4779
$ws.Range("L129").Value = 8662.667099999999
$ws.Range("M129").Value = 221
$ws.Range("N129").Value = -18662.6671
$ws.Range("H130").Value = 3463.4285
$ws.Range("I130").Value = 936.25
$ws.Range("J130").Value = 6833
$ws.Range("K130").Value = 2808.75
$ws.Range("L130").Value = 20499
$ws.Range("M130").Value = 2211.25
$ws.Range("N130").Value = -30539
$ws.Range("H131").Value = 1072.5
$ws.Range("J131").Value = 1553
$ws.Range("L131").Value = 4659
$ws.Range("N131").Value = -14739
$ws.Range("H134").Value = 1907.5454
$ws.Range("I134").Value = 1907.5454
$ws.Range("K134").Value = 5722.6362
$ws.Range("M134").Value = -652.6361999999999
$ws.Range("H136").Value = 11044.556
$ws.Range("I136").Value = 6566.8335
$ws.Range("K136").Value = 19700.5005
$ws.Range("M136").Value = -14600.5005
$ws.Range("H137").Value = 6166
$ws.Range("I137").Value = 4952.8
$ws.Range("J137").Value = 7379.2
$ws.Range("K137").Value = 14858.4
$ws.Range("L137").Value = 22137.6
$ws.Range("M137").Value = -9758.400000000001
$ws.Range("N137").Value = -32337.6
$ws.Range("H138").Value = 4249.5
$ws.Range("I138").Value = 4249.5
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 12748.5
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -7608.5
$ws.Range("N138").Value = ""
$ws.Range("H139").Value = 2467.5454
$ws.Range("I139").Value = 2467.5454
$ws.Range("K139").Value = 7402.6362
$ws.Range("M139").Value = -2262.6362
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2812.5642
$ws.Range("I132").Value = 2716.5881
$ws.Range("K132").Value = 8149.7643
$ws.Range("M132").Value = -5619.7643
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 4004398.5
$ws.Range("I2").Value = 10006500
$ws.Range("J2").Value = 2997.6667
$ws.Range("K2").Value = 10006500
$ws.Range("L2").Value = 2997.6667
$ws.Range("M2").Value = -10006388
$ws.Range("N2").Value = -3221.6667
$ws.Range("H22").Value = 8020.92
$ws.Range("I22").Value = 5752.231
$ws.Range("J22").Value = 10478.667
$ws.Range("K22").Value = 5752.231
$ws.Range("L22").Value = 10478.667
$ws.Range("M22").Value = -5457.231
$ws.Range("N22").Value = -11068.667
$ws.Range("H27").Value = 8020.92
$ws.Range("I27").Value = 5752.231
$ws.Range("J27").Value = 10478.667
$ws.Range("K27").Value = 5752.231
$ws.Range("L27").Value = 10478.667
$ws.Range("M27").Value = -5645.231
$ws.Range("N27").Value = -10692.667
$ws.Range("H46").Value = 2737.75
$ws.Range("I46").Value = 2143.5
$ws.Range("K46").Value = 2143.5
$ws.Range("M46").Value = -1955.5
$ws.Range("H132").Value = 5786.125
$ws.Range("I132").Value = 4904.7144
$ws.Range("K132").Value = 14714.1432
$ws.Range("M132").Value = -12184.1432
$ws.Range("H136").Value = 4697
$ws.Range("I136").Value = 3500
$ws.Range("J136").Value = 5894
$ws.Range("K136").Value = 10500
$ws.Range("L136").Value = 17682
$ws.Range("M136").Value = -7950
$ws.Range("N136").Value = -22782
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1537.6471
$ws.Range("I107").Value = 1426.4166
$ws.Range("J107").Value = 1804.6
$ws.Range("K107").Value = 4279.2498
$ws.Range("L107").Value = 5413.799999999999
$ws.Range("M107").Value = -2359.2498
$ws.Range("N107").Value = -9253.799999999999
$ws.Range("H122").Value = 983.44446
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").Value = ""
